$wb = $excel.ActiveWorkbook

# --- Sheet: All Estimates ---
$ws = $wb.Worksheets.Item("All Estimates")
$ws.Range("D2").Value = -0.45
$ws.Range("E2").Value = -0.78
$ws.Range("F2").Value = -0.11
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("A3").Value = "GLM"
$ws.Range("C3").Value = "Indirect"
$ws.Range("D3").Value = -0.01
$ws.Range("E3").Value = -0.48
$ws.Range("F3").Value = 0.46
$ws.Range("G3").Value = 0.97
$ws.Range("H3").Value = 0.06
$ws.Range("I3").Value = 0.24
$ws.Range("J3").ClearContents()
$ws.Range("A4").Value = "GLM"
$ws.Range("C4").Value = "Indirect Calculated"
$ws.Range("D4").Value = -0.01
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("C5").Value = "Total"
$ws.Range("D5").Value = -0.46
$ws.Range("E5").Value = -0.79
$ws.Range("F5").Value = -0.12
$ws.Range("G5").Value = 0.01
$ws.Range("H5").Value = 0.03
$ws.Range("I5").Value = 0.17
$ws.Range("A6").Value = "GLM"
$ws.Range("C6").Value = "PM"
$ws.Range("D6").Value = 0.02
$ws.Range("E6").Value = -1
$ws.Range("F6").Value = 1.05
$ws.Range("G6").Value = 0.48
$ws.Range("H6").Value = 0.27
$ws.Range("I6").Value = 0.52
$ws.Range("A7").Value = "GLM"
$ws.Range("C7").Value = "PM Calculated"
$ws.Range("D7").Value = 0.02
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("A8").Value = "GLMM"
$ws.Range("C8").Value = "Direct"
$ws.Range("D8").Value = -0.45
$ws.Range("E8").Value = -0.83
$ws.Range("F8").Value = -0.07000000000000001
$ws.Range("G8").Value = 0.02
$ws.Range("H8").Value = 0.04
$ws.Range("I8").Value = 0.19
$ws.Range("J8").Value = 0.02
$ws.Range("C9").Value = "Indirect"
$ws.Range("E9").Value = -0.55
$ws.Range("F9").Value = 0.53
$ws.Range("G9").Value = 0.97
$ws.Range("H9").Value = 0.07000000000000001
$ws.Range("I9").Value = 0.27
$ws.Range("A10").Value = "GLMM"
$ws.Range("D10").Value = -0.01
$ws.Range("A11").Value = "GLMM"
$ws.Range("D11").Value = -0.46
$ws.Range("E11").Value = -0.84
$ws.Range("F11").Value = -0.08
$ws.Range("G11").Value = 0.02
$ws.Range("H11").Value = 0.04
$ws.Range("I11").Value = 0.19
$ws.Range("J11").Value = 0.02
$ws.Range("C12").Value = "PM"
$ws.Range("D12").Value = 0.02
$ws.Range("E12").Value = -1.13
$ws.Range("F12").Value = 1.18
$ws.Range("G12").Value = 0.49
$ws.Range("H12").Value = 0.35
$ws.Range("I12").Value = 0.59
$ws.Range("J12").ClearContents()
$ws.Range("A13").Value = "GLMM"
$ws.Range("C13").Value = "PM Calculated"
$ws.Range("D13").Value = 0.02
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("A14").Value = "GEE"
$ws.Range("C14").Value = "Direct"
$ws.Range("D14").Value = -0.45
$ws.Range("E14").Value = -0.8
$ws.Range("F14").Value = -0.09
$ws.Range("G14").Value = 0.01
$ws.Range("H14").Value = 0.03
$ws.Range("I14").Value = 0.18
$ws.Range("J14").Value = 0
$ws.Range("A15").Value = "GEE"
$ws.Range("C15").Value = "Indirect"
$ws.Range("D15").Value = -0.01
$ws.Range("E15").Value = -0.52
$ws.Range("F15").Value = 0.49
$ws.Range("G15").Value = 0.96
$ws.Range("H15").Value = 0.07000000000000001
$ws.Range("I15").Value = 0.26
$ws.Range("C16").Value = "Indirect Calculated"
$ws.Range("D16").Value = -0.01
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("A17").Value = "GEE"
$ws.Range("C17").Value = "Total"
$ws.Range("D17").Value = -0.46
$ws.Range("E17").Value = -0.82
$ws.Range("F17").Value = -0.11
$ws.Range("G17").Value = 0.01
$ws.Range("H17").Value = 0.03
$ws.Range("I17").Value = 0.18
$ws.Range("J17").Value = 0
$ws.Range("A18").Value = "GEE"
$ws.Range("C18").Value = "PM"
$ws.Range("D18").Value = 0.03
$ws.Range("E18").Value = -1.05
$ws.Range("F18").Value = 1.1
$ws.Range("G18").Value = 0.48
$ws.Range("H18").Value = 0.3
$ws.Range("I18").Value = 0.55
$ws.Range("D19").Value = 0.03

# --- Sheet: All Estimates (OR) ---
$ws = $wb.Worksheets.Item("All Estimates (OR)")
$ws.Range("E2").Value = 0.46
$ws.Range("F2").Value = 0.89
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("A3").Value = "GLM"
$ws.Range("C3").Value = "Indirect"
$ws.Range("D3").Value = 0.99
$ws.Range("E3").Value = 0.62
$ws.Range("F3").Value = 1.59
$ws.Range("G3").Value = 0.97
$ws.Range("H3").Value = 0.06
$ws.Range("I3").Value = 0.24
$ws.Range("J3").ClearContents()
$ws.Range("A4").Value = "GLM"
$ws.Range("C4").Value = "Indirect Calculated"
$ws.Range("D4").Value = 0.99
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("C5").Value = "Total"
$ws.Range("D5").Value = 0.63
$ws.Range("E5").Value = 0.45
$ws.Range("F5").Value = 0.88
$ws.Range("G5").Value = 0.01
$ws.Range("H5").Value = 0.03
$ws.Range("I5").Value = 0.17
$ws.Range("A6").Value = "GLM"
$ws.Range("C6").Value = "PM"
$ws.Range("D6").Value = 0.02
$ws.Range("E6").Value = -1
$ws.Range("F6").Value = 1.05
$ws.Range("G6").Value = 0.48
$ws.Range("H6").Value = 0.27
$ws.Range("I6").Value = 0.52
$ws.Range("A7").Value = "GLM"
$ws.Range("C7").Value = "PM Calculated"
$ws.Range("D7").Value = 0.02
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("A8").Value = "GLMM"
$ws.Range("C8").Value = "Direct"
$ws.Range("D8").Value = 0.64
$ws.Range("E8").Value = 0.44
$ws.Range("F8").Value = 0.93
$ws.Range("G8").Value = 0.02
$ws.Range("H8").Value = 0.04
$ws.Range("I8").Value = 0.19
$ws.Range("J8").Value = 0.02
$ws.Range("C9").Value = "Indirect"
$ws.Range("E9").Value = 0.58
$ws.Range("F9").Value = 1.69
$ws.Range("G9").Value = 0.97
$ws.Range("H9").Value = 0.07000000000000001
$ws.Range("I9").Value = 0.27
$ws.Range("A10").Value = "GLMM"
$ws.Range("D10").Value = 0.99
$ws.Range("A11").Value = "GLMM"
$ws.Range("D11").Value = 0.63
$ws.Range("E11").Value = 0.43
$ws.Range("F11").Value = 0.92
$ws.Range("G11").Value = 0.02
$ws.Range("H11").Value = 0.04
$ws.Range("I11").Value = 0.19
$ws.Range("J11").Value = 0.02
$ws.Range("C12").Value = "PM"
$ws.Range("D12").Value = 0.02
$ws.Range("E12").Value = -1.13
$ws.Range("F12").Value = 1.18
$ws.Range("G12").Value = 0.49
$ws.Range("H12").Value = 0.35
$ws.Range("I12").Value = 0.59
$ws.Range("J12").ClearContents()
$ws.Range("A13").Value = "GLMM"
$ws.Range("C13").Value = "PM Calculated"
$ws.Range("D13").Value = 0.02
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()
$ws.Range("I13").ClearContents()
$ws.Range("J13").ClearContents()
$ws.Range("A14").Value = "GEE"
$ws.Range("C14").Value = "Direct"
$ws.Range("D14").Value = 0.64
$ws.Range("E14").Value = 0.45
$ws.Range("F14").Value = 0.91
$ws.Range("G14").Value = 0.01
$ws.Range("H14").Value = 0.03
$ws.Range("I14").Value = 0.18
$ws.Range("J14").Value = 0
$ws.Range("A15").Value = "GEE"
$ws.Range("C15").Value = "Indirect"
$ws.Range("D15").Value = 0.99
$ws.Range("E15").Value = 0.6
$ws.Range("F15").Value = 1.63
$ws.Range("G15").Value = 0.96
$ws.Range("H15").Value = 0.07000000000000001
$ws.Range("I15").Value = 0.26
$ws.Range("C16").Value = "Indirect Calculated"
$ws.Range("D16").Value = 0.99
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("A17").Value = "GEE"
$ws.Range("C17").Value = "Total"
$ws.Range("D17").Value = 0.63
$ws.Range("E17").Value = 0.44
$ws.Range("F17").Value = 0.9
$ws.Range("G17").Value = 0.01
$ws.Range("H17").Value = 0.03
$ws.Range("I17").Value = 0.18
$ws.Range("J17").Value = 0
$ws.Range("A18").Value = "GEE"
$ws.Range("C18").Value = "PM"
$ws.Range("D18").Value = 0.03
$ws.Range("E18").Value = -1.05
$ws.Range("F18").Value = 1.1
$ws.Range("G18").Value = 0.48
$ws.Range("H18").Value = 0.3
$ws.Range("I18").Value = 0.55
$ws.Range("D19").Value = 0.03

# --- Sheet: Total Effects ---
$ws = $wb.Worksheets.Item("Total Effects")
$ws.Range("D2").Value = -0.46
$ws.Range("E2").Value = -0.79
$ws.Range("F2").Value = -0.12
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("D3").Value = -0.46
$ws.Range("E3").Value = -0.84
$ws.Range("F3").Value = -0.08
$ws.Range("G3").Value = 0.02
$ws.Range("H3").Value = 0.04
$ws.Range("I3").Value = 0.19
$ws.Range("J3").Value = 0.02
$ws.Range("D4").Value = -0.46
$ws.Range("E4").Value = -0.82
$ws.Range("F4").Value = -0.11
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 0.03
$ws.Range("I4").Value = 0.18
$ws.Range("J4").Value = 0

# --- Sheet: Total Effects (OR) ---
$ws = $wb.Worksheets.Item("Total Effects (OR)")
$ws.Range("D2").Value = 0.63
$ws.Range("E2").Value = 0.45
$ws.Range("F2").Value = 0.88
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("D3").Value = 0.63
$ws.Range("E3").Value = 0.43
$ws.Range("F3").Value = 0.92
$ws.Range("G3").Value = 0.02
$ws.Range("H3").Value = 0.04
$ws.Range("I3").Value = 0.19
$ws.Range("J3").Value = 0.02
$ws.Range("E4").Value = 0.44
$ws.Range("F4").Value = 0.9
$ws.Range("G4").Value = 0.01
$ws.Range("H4").Value = 0.03
$ws.Range("I4").Value = 0.18
$ws.Range("J4").Value = 0

# --- Sheet: Direct Effects ---
$ws = $wb.Worksheets.Item("Direct Effects")
$ws.Range("D2").Value = -0.45
$ws.Range("E2").Value = -0.78
$ws.Range("F2").Value = -0.11
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("D3").Value = 0.58
$ws.Range("E3").Value = 0.22
$ws.Range("F3").Value = 0.9399999999999999
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.03
$ws.Range("I3").Value = 0.18
$ws.Range("D4").Value = -0.45
$ws.Range("E4").Value = -0.83
$ws.Range("H4").Value = 0.04
$ws.Range("I4").Value = 0.19
$ws.Range("J4").Value = 0.02
$ws.Range("D5").Value = 0.57
$ws.Range("E5").Value = 0.21
$ws.Range("F5").Value = 0.93
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.03
$ws.Range("I5").Value = 0.18
$ws.Range("J5").Value = 0.02
$ws.Range("D6").Value = -0.45
$ws.Range("E6").Value = -0.8
$ws.Range("F6").Value = -0.09
$ws.Range("G6").Value = 0.01
$ws.Range("H6").Value = 0.03
$ws.Range("I6").Value = 0.18
$ws.Range("J6").Value = 0
$ws.Range("D7").Value = 0.57
$ws.Range("E7").Value = 0.24
$ws.Range("F7").Value = 0.9
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0.03
$ws.Range("I7").Value = 0.17
$ws.Range("J7").Value = 0

# --- Sheet: Direct Effects (OR) ---
$ws = $wb.Worksheets.Item("Direct Effects (OR)")
$ws.Range("E2").Value = 0.46
$ws.Range("F2").Value = 0.89
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.03
$ws.Range("I2").Value = 0.17
$ws.Range("D3").Value = 1.78
$ws.Range("E3").Value = 1.24
$ws.Range("F3").Value = 2.55
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.03
$ws.Range("I3").Value = 0.18
$ws.Range("D4").Value = 0.64
$ws.Range("E4").Value = 0.44
$ws.Range("F4").Value = 0.93
$ws.Range("H4").Value = 0.04
$ws.Range("I4").Value = 0.19
$ws.Range("J4").Value = 0.02
$ws.Range("D5").Value = 1.77
$ws.Range("E5").Value = 1.24
$ws.Range("F5").Value = 2.54
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0.03
$ws.Range("I5").Value = 0.18
$ws.Range("J5").Value = 0.02
$ws.Range("D6").Value = 0.64
$ws.Range("E6").Value = 0.45
$ws.Range("G6").Value = 0.01
$ws.Range("H6").Value = 0.03
$ws.Range("I6").Value = 0.18
$ws.Range("J6").Value = 0
$ws.Range("D7").Value = 1.77
$ws.Range("E7").Value = 1.27
$ws.Range("F7").Value = 2.47
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0.03
$ws.Range("I7").Value = 0.17
$ws.Range("J7").Value = 0

